$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 48.42857
$ws.Range("I5").Value = 36.666668
$ws.Range("J5").Value = 119
$ws.Range("K5").Value = 36.666668
$ws.Range("L5").Value = 119
$ws.Range("M5").Value = 78.333332
$ws.Range("N5").Value = -349
$ws.Range("I33").Value = 2103.9092
$ws.Range("J33").Value = 1371
$ws.Range("K33").Value = 2103.9092
$ws.Range("L33").Value = 1371
$ws.Range("M33").Value = -1874.9092
$ws.Range("N33").Value = -1829
$ws.Range("H40").Value = 5143
$ws.Range("I40").Value = 4029.8
$ws.Range("J40").Value = 6998.3335
$ws.Range("K40").Value = 4029.8
$ws.Range("L40").Value = 6998.3335
$ws.Range("M40").Value = -3854.8
$ws.Range("N40").Value = -7348.3335
$ws.Range("H43").Value = 201140.42
$ws.Range("I43").Value = 6525
$ws.Range("K43").Value = 6525
$ws.Range("M43").Value = -6456
$ws.Range("H58").Value = 57696572
$ws.Range("J58").Value = 62506804
$ws.Range("L58").Value = 187520412
$ws.Range("N58").Value = -187520712
$ws.Range("H75").Value = 41657
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 41657
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H80").Value = 4049442.2
$ws.Range("I80").Value = 6584913
$ws.Range("K80").Value = 19754739
$ws.Range("M80").Value = -19753741
$ws.Range("H83").Value = 4049442.2
$ws.Range("I83").Value = 6584913
$ws.Range("K83").Value = 59264217
$ws.Range("M83").Value = -59259225
$ws.Range("H113").Value = 146485200
$ws.Range("I113").Value = 111111110
$ws.Range("J113").Value = 150022610
$ws.Range("K113").Value = 111111110
$ws.Range("L113").Value = 150022610
$ws.Range("M113").Value = -111107856
$ws.Range("N113").Value = -150029118
$ws.Range("H115").Value = 2538.111
$ws.Range("I115").Value = 281.33334
$ws.Range("J115").Value = 3666.5
$ws.Range("K115").Value = 844.0000200000001
$ws.Range("L115").Value = 10999.5
$ws.Range("M115").Value = 722.9999799999999
$ws.Range("N115").Value = -14133.5
$ws.Range("H116").Value = 17866622
$ws.Range("I116").Value = 83336770
$ws.Range("K116").Value = 83336770
$ws.Range("M116").Value = -83333328
$ws.Range("H132").Value = 1430.9807
$ws.Range("I132").Value = 1295.9788
$ws.Range("J132").Value = 2700
$ws.Range("K132").Value = 3887.936400000001
$ws.Range("L132").Value = 8100
$ws.Range("M132").Value = -1357.936400000001
$ws.Range("N132").Value = -13160
$ws.Range("H137").Value = 3508.0513
$ws.Range("I137").Value = 3145.6
$ws.Range("K137").Value = 9436.799999999999
$ws.Range("M137").Value = -6886.799999999999
$ws.Range("H138").Value = 1566435.9
$ws.Range("I138").Value = 1060
$ws.Range("J138").Value = 3036940.5
$ws.Range("K138").Value = 3180
$ws.Range("L138").Value = 9110821.5
$ws.Range("M138").Value = 1960
$ws.Range("N138").Value = -9121101.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1989946.8
$ws.Range("I32").Value = 2159342
$ws.Range("K32").Value = 2159342
$ws.Range("M32").Value = -2159055
$ws.Range("H45").Value = 4888.9
$ws.Range("I45").Value = 2646
$ws.Range("K45").Value = 2646
$ws.Range("M45").Value = -2269
$ws.Range("H74").Value = 23204.098
$ws.Range("I74").Value = 31969
$ws.Range("J74").Value = 4797.8
$ws.Range("K74").Value = 31969
$ws.Range("L74").Value = 4797.8
$ws.Range("M74").Value = -31095
$ws.Range("N74").Value = -6545.8
$ws.Range("H77").Value = 23204.098
$ws.Range("I77").Value = 31969
$ws.Range("J77").Value = 4797.8
$ws.Range("K77").Value = 159845
$ws.Range("L77").Value = 23989
$ws.Range("M77").Value = -155477
$ws.Range("N77").Value = -32725
$ws.Range("H108").Value = 52188
$ws.Range("J108").Value = 52188
$ws.Range("L108").Value = 52188
$ws.Range("N108").Value = -59868
$ws.Range("H132").Value = 4387.6
$ws.Range("I132").Value = 1877.8206
$ws.Range("K132").Value = 5633.4618
$ws.Range("M132").Value = -3103.4618

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 273.33334
$ws.Range("I22").Value = 273.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 273.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -100.33334
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 1525.0322
$ws.Range("I94").Value = 1266.92
$ws.Range("K94").Value = 1266.92
$ws.Range("M94").Value = -815.9200000000001
$ws.Range("H105").Value = 3310.2856
$ws.Range("I105").Value = 2580.5715
$ws.Range("J105").Value = 4769.7144
$ws.Range("K105").Value = 2580.5715
$ws.Range("L105").Value = 4769.7144
$ws.Range("M105").Value = -833.5715
$ws.Range("N105").Value = -8263.714400000001
$ws.Range("H113").Value = 5268.25
$ws.Range("I113").Value = 5268.25
$ws.Range("K113").Value = 5268.25
$ws.Range("M113").Value = -3098.25
$ws.Range("H134").Value = 4303.9277
$ws.Range("J134").Value = 11141
$ws.Range("L134").Value = 33423
$ws.Range("N134").Value = -38493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 103.83871
$ws.Range("I7").Value = 54.761906
$ws.Range("K7").Value = 54.761906
$ws.Range("M7").Value = 58.238094
$ws.Range("H15").Value = 15573.8
$ws.Range("I15").Value = 956.3333
$ws.Range("J15").Value = 37500
$ws.Range("K15").Value = 956.3333
$ws.Range("L15").Value = 37500
$ws.Range("M15").Value = -786.3333
$ws.Range("N15").Value = -37840
$ws.Range("H31").Value = 5561912.5
$ws.Range("I31").Value = 2744.2903
$ws.Range("J31").Value = 11504472
$ws.Range("K31").Value = 2744.2903
$ws.Range("L31").Value = 11504472
$ws.Range("M31").Value = -2449.2903
$ws.Range("N31").Value = -11505062
$ws.Range("H34").Value = 5561912.5
$ws.Range("I34").Value = 2744.2903
$ws.Range("J34").Value = 11504472
$ws.Range("K34").Value = 2744.2903
$ws.Range("L34").Value = 11504472
$ws.Range("M34").Value = -2542.2903
$ws.Range("N34").Value = -11504876
$ws.Range("H132").Value = 9097395
$ws.Range("I132").Value = 4045.1428
$ws.Range("K132").Value = 12135.4284
$ws.Range("M132").Value = -9605.428400000001
$ws.Range("H134").Value = 4168.039
$ws.Range("I134").Value = 1711.0435
$ws.Range("K134").Value = 5133.1305
$ws.Range("M134").Value = -2598.1305

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 76230504
$ws.Range("J4").Value = 2171477
$ws.Range("L4").Value = 6514431
$ws.Range("N4").Value = -6514655
$ws.Range("H110").Value = 461
$ws.Range("I110").Value = 461
$ws.Range("K110").Value = 1383
$ws.Range("M110").Value = 2707
$ws.Range("H132").Value = 10339.143
$ws.Range("J132").Value = 18566.416
$ws.Range("L132").Value = 167097.744
$ws.Range("N132").Value = -172157.744
$ws.Range("H136").Value = 1575.1111
$ws.Range("I136").Value = 1575.1111
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4725.3333
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 374.6666999999998
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 6719.5
$ws.Range("I141").Value = 2959.3333
$ws.Range("K141").Value = 8877.999899999999
$ws.Range("M141").Value = -3697.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 11202601
$ws.Range("I14").Value = 28000502
$ws.Range("J14").Value = 4000
$ws.Range("K14").Value = 28000502
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = -28000334
$ws.Range("N14").Value = -4336
$ws.Range("H80").Value = 128279.875
$ws.Range("I80").Value = 3966.5
$ws.Range("K80").Value = 3966.5
$ws.Range("M80").Value = -2968.5
$ws.Range("H83").Value = 128279.875
$ws.Range("I83").Value = 3966.5
$ws.Range("K83").Value = 19832.5
$ws.Range("M83").Value = -14840.5
$ws.Range("H94").Value = 30960.375
$ws.Range("J94").Value = 30960.375
$ws.Range("L94").Value = 30960.375
$ws.Range("N94").Value = -32312.375
$ws.Range("H97").Value = 1753.4286
$ws.Range("I97").Value = 1315.3334
$ws.Range("J97").Value = 2542
$ws.Range("K97").Value = 1315.3334
$ws.Range("L97").Value = 2542
$ws.Range("M97").Value = -819.3334
$ws.Range("N97").Value = -3534

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 1408
$ws.Range("I17").Value = 1408
$ws.Range("K17").Value = 1408
$ws.Range("M17").Value = -1238
$ws.Range("H46").Value = 4834209.5
$ws.Range("I46").Value = 1733.2222
$ws.Range("K46").Value = 1733.2222
$ws.Range("M46").Value = -1545.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 10753565
$ws.Range("J107").Value = 25642108
$ws.Range("L107").Value = 76926324
$ws.Range("N107").Value = -76930164
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 12199477
$ws.Range("I132").Value = 14710051
$ws.Range("K132").Value = 44130153
$ws.Range("M132").Value = -44127623
$ws.Range("H136").Value = 19631026
$ws.Range("I136").Value = 33334112
$ws.Range("K136").Value = 100002336
$ws.Range("M136").Value = -99999786
